# Add "amountVND" (number written in words) + "Kế Toán" signature block
# to the "Phieu De Nghi Thanh Toan" (refund/advance request) template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 17 ("Ban Giám Đốc" / approval block),
# pushing the existing rows 17-18 down to 18-19.
$ws.Rows("17:17").Insert()

# Populate the newly freed row 17 with the "amount in words" label/value pair,
# mirroring the style of the "grand total" row right above it (row 16).
$ws.Range("C17").Value2 = "Số tiền bằng chữ:"
$ws.Range("E17").Value2 = '${amountVND}'

# Row 19 (previously row 18) currently reads:
#   B19=Phê duyệt | C19=Trưởng Bộ phận | D19=Người Đề nghị
# We need to insert a new "Kế Toán" signer between "Trưởng Bộ phận" and
# "Người Đề nghị", so first copy the existing D19 cell (value + style) into
# the new E19 slot, then overwrite D19 with the new "Kế Toán" label.
$ws.Range("D19").Copy($ws.Range("E19"))
$ws.Range("D19").Value2 = "Kế Toán"

# Reflect the new selection left behind by the edit (active cell = C17).
$ws.Range("C17").Select() | Out-Null
